$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 74; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $cell.Value2 + 1
}
